# Update "想去人数" (F column) counts on each sheet per the latest scrape.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 217
$ws1.Range("F3").Value = 702
$ws1.Range("F5").Value = 2226
$ws1.Range("F6").Value = 1319
$ws1.Range("F8").Value = 90
$ws1.Range("F10").Value = 2855
$ws1.Range("F13").Value = 1071
$ws1.Range("F14").Value = 574
$ws1.Range("F16").Value = 912
$ws1.Range("F17").Value = 94
$ws1.Range("F20").Value = 124
$ws1.Range("F22").Value = 584
$ws1.Range("F23").Value = 275
$ws1.Range("F25").Value = 971
$ws1.Range("F26").Value = 4882
$ws1.Range("F27").Value = 389
$ws1.Range("F28").Value = 150
$ws1.Range("F29").Value = 70
$ws1.Range("F30").Value = 91

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 386
$ws2.Range("F11").Value = 183
$ws2.Range("F20").Value = 28
$ws2.Range("F24").Value = 349
$ws2.Range("F26").Value = 566
$ws2.Range("F36").Value = 712
$ws2.Range("F37").Value = 32

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 381
$ws3.Range("F7").Value = 349

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 217
$ws4.Range("F6").Value = 381
$ws4.Range("F8").Value = 702
$ws4.Range("F10").Value = 386
$ws4.Range("F13").Value = 2226
$ws4.Range("F14").Value = 1319
$ws4.Range("F16").Value = 90
$ws4.Range("F18").Value = 183
$ws4.Range("F20").Value = 2855
$ws4.Range("F24").Value = 1071
$ws4.Range("F25").Value = 574
$ws4.Range("F27").Value = 349
$ws4.Range("F29").Value = 912
$ws4.Range("F30").Value = 912
$ws4.Range("F31").Value = 94
$ws4.Range("F32").Value = 28
$ws4.Range("F35").Value = 124
$ws4.Range("F39").Value = 584
$ws4.Range("F40").Value = 349
$ws4.Range("F41").Value = 275
$ws4.Range("F44").Value = 971
$ws4.Range("F45").Value = 4882
$ws4.Range("F47").Value = 389
$ws4.Range("F48").Value = 150
$ws4.Range("F49").Value = 712
$ws4.Range("F50").Value = 91
